$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet holds a daily log table in A1:D39 (date, day-of-week, hour,
# ranking). A new day's entry needs to be appended as row 40:
#   2025/09/30, 火, 16, 3
# Column A holds dates stored as literal text (e.g. "2025/09/30"), so we
# copy the existing A39 cell (same text) down into A40 instead of typing
# the string directly - that avoids Excel's automatic text-to-date
# conversion while keeping the cell's formatting untouched.
$ws.Range("A39").Copy($ws.Range("A40"))

$ws.Cells.Item(40, 2).Value = "火"
$ws.Cells.Item(40, 3).Value = 16
$ws.Cells.Item(40, 4).Value = 3
